# json schema und daten cleanup
# Adds a 10th row of VIP-night pairing data (mirroring the existing
# rotation pattern) to the "Nights" sheet, and updates the sheet's
# selection/scroll state to match.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 10 values, following the same rotation seen in rows 6-9
# (green-highlighted group).
$ws.Range("B10").Value = "Gabriela"
$ws.Range("C10").Value = "Asena"
$ws.Range("D10").Value = "Emmy"
$ws.Range("E10").Value = "Nadja"
$ws.Range("F10").Value = "Jennifer"
$ws.Range("G10").Value = "Tara"
$ws.Range("H10").Value = "Laura M."
$ws.Range("I10").Value = "Laura L."
$ws.Range("J10").Value = "Anastasia"
$ws.Range("K10").Value = "Linda"
$ws.Range("L10").Value = 10

# Match the formatting already used for the row-6..9 block: green fill
# across B:K, and the right-aligned numeric style in column L.
$ws.Range("B10:K10").Interior.Color = 5287936
$ws.Range("L10").NumberFormat = "#,##0"
$ws.Range("L10").HorizontalAlignment = -4152

# Update the view: selection moves to the newly added row, and the
# sheet scrolls back so column A is visible again (topLeftCell cleared).
$ws.Range("B9:B10").Select() | Out-Null
